# Applies the update described by the commit:
#  - Swap the match-detail columns (F:V) between several row pairs
#    (these rows keep their Indice/pais/torneio/temporada/data_partida
#    columns A:E - only the match content in F:V moved rows).
#  - Append a new match row (144) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V content between row pairs -----------------------------------
$pairs = @(
    @(7, 8),
    @(34, 36),
    @(35, 37),
    @(39, 40),
    @(41, 42),
    @(56, 57)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]
    $rng1 = $ws.Range("F$row1" + ":V$row1")
    $rng2 = $ws.Range("F$row2" + ":V$row2")
    $val1 = $rng1.Value2
    $val2 = $rng2.Value2
    $rng1.Value2 = $val2
    $rng2.Value2 = $val1
}

# --- Append the new row (144) with the latest scraped match ---------------
# Column A uses the same bold/bordered/centered style as every other row's
# "Indice" cell, and column E uses the date-time number format - grab both
# by copying the formatting from the previous row before writing the values
# (avoids minting new, divergent style entries).
$newRow = 144
$prevRow = 143

$ws.Cells.Item($prevRow, 1).Copy($ws.Cells.Item($newRow, 1))
$ws.Cells.Item($newRow, 1).Value2 = 143

$ws.Cells.Item($newRow, 2).Value = "poland"
$ws.Cells.Item($newRow, 3).Value = "division-2"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"

$ws.Cells.Item($prevRow, 5).Copy($ws.Cells.Item($newRow, 5))
$ws.Cells.Item($newRow, 5).Value2 = 45240.75

$ws.Cells.Item($newRow, 6).Value = "Kotwica Kolobrzeg"
$ws.Cells.Item($newRow, 7).Value2 = 2
$ws.Cells.Item($newRow, 8).Value = "Sandecja Nowy S."
$ws.Cells.Item($newRow, 9).Value2 = 1
$ws.Cells.Item($newRow, 10).Value2 = 1.6
$ws.Cells.Item($newRow, 11).Value = "09/11/2023 06:12"
$ws.Cells.Item($newRow, 12).Value2 = 1.56
$ws.Cells.Item($newRow, 13).Value = "10/11/2023 17:54"
$ws.Cells.Item($newRow, 14).Value2 = 3.77
$ws.Cells.Item($newRow, 15).Value = "09/11/2023 06:12"
$ws.Cells.Item($newRow, 16).Value2 = 3.99
$ws.Cells.Item($newRow, 17).Value = "10/11/2023 17:57"
$ws.Cells.Item($newRow, 18).Value2 = 4.45
$ws.Cells.Item($newRow, 19).Value = "09/11/2023 06:12"
$ws.Cells.Item($newRow, 20).Value2 = 5.58
$ws.Cells.Item($newRow, 21).Value = "10/11/2023 17:55"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/poland/division-2/kotwica-kolobrzeg-sandecja-nowy-s/MLu2BTl2/"
